# rmarkdown.docx -- "lectures for day 6"
#
# The rendered R values in the "Inline Code" example were re-knit, so the
# literal numbers baked into the prose need to be refreshed to match the
# newly-rendered output.

$d = $word.ActiveDocument

$old = "Then I can write x = -0.7336645 and y = 0.9959624."
$new = "Then I can write x = 2.1614889 and y = 1.3853053."

$found = $d.Content.Find.Execute(
    $old,   # FindText
    $true,  # MatchCase
    $true,  # MatchWholeWord
    $false, # MatchWildcards
    $false, # MatchSoundsLike
    $false, # MatchAllWordForms
    $true,  # Forward
    1,      # Wrap (wdFindContinue)
    $false, # Format
    $new,   # ReplaceWith
    2       # Replace (wdReplaceAll)
)

Write-Output "Replaced x/y values: $found"
